# Add data for 2022-04-24 — bump the "through April 15" window to "through April 16"
# and add the one day's worth of newly-recorded carjackings to the matching cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the column header label for the current month.
$ws.Name = "Through 2022-04-16"
$ws.Range("B1").Value = "April 2022 (through April 16)"

# Update/insert cell values in column B (the "April 2022" column) and the
# handful of other neighborhood/month cells touched by the day's new records.
$ws.Range("N2").Value = 4
$ws.Range("V2").Value = 2
$ws.Range("B3").Value = 4
$ws.Range("F4").Value = 4
$ws.Range("R4").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("J6").Value = 1
$ws.Range("V11").Value = 1
$ws.Range("B17").Value = 2
$ws.Range("B32").Value = 1
$ws.Range("F32").Value = 3
$ws.Range("Z32").Value = 2
$ws.Range("R48").Value = 1
$ws.Range("B51").Value = 2
$ws.Range("V60").Value = 1
$ws.Range("F66").Value = 1
$ws.Range("B75").Value = 2
